$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.779501060154568
$ws.Range("D2").Value = 10.2373938384862
$ws.Range("E2").Value = 13.99612960964306
$ws.Range("F2").Value = 27.22541788187684
$ws.Range("G2").Value = 26.44254053043482
$ws.Range("H2").Value = 13.15773181025866
$ws.Range("J2").Value = 9.82979599056854
$ws.Range("M2").Value = 18.27105920890671
$ws.Range("N2").Value = 17.83519077232751
$ws.Range("O2").Value = 19.81460210027505
$ws.Range("B3").Value = 7.705586864444626
$ws.Range("D3").Value = 10.26968971567905
$ws.Range("E3").Value = 14.06525292659726
$ws.Range("F3").Value = 27.1279042206101
$ws.Range("G3").Value = 26.0463120115921
$ws.Range("H3").Value = 13.15993155945395
$ws.Range("J3").Value = 9.870375185156
$ws.Range("M3").Value = 17.75885792108064
$ws.Range("N3").Value = 17.74420608268311
$ws.Range("O3").Value = 19.74533093088947
$ws.Range("B4").Value = 7.661628248227161
$ws.Range("D4").Value = 10.29108972175371
$ws.Range("E4").Value = 14.11024113680611
$ws.Range("F4").Value = 27.07559584471297
$ws.Range("G4").Value = 25.80933067477273
$ws.Range("H4").Value = 13.16395093012061
$ws.Range("J4").Value = 9.896594097422296
$ws.Range("M4").Value = 17.43702975154642
$ws.Range("N4").Value = 17.69051490325329
$ws.Range("O4").Value = 19.70792797885207
$ws.Range("B5").Value = 7.64409362305902
$ws.Range("D5").Value = 10.30020523727466
$ws.Range("E5").Value = 14.12921480063321
$ws.Range("F5").Value = 27.05619828320707
$ws.Range("G5").Value = 25.71448120192073
$ws.Range("H5").Value = 13.166260293027
$ws.Range("J5").Value = 9.907606917723378
$ws.Range("M5").Value = 17.30420987387504
$ws.Range("N5").Value = 17.6691998329978
$ws.Range("O5").Value = 19.69398725434686
$ws.Range("B6").Value = 7.641205481415704
$ws.Range("D6").Value = 10.30174270898357
$ws.Range("E6").Value = 14.13240406433682
$ws.Range("F6").Value = 27.053093628025
$ws.Range("G6").Value = 25.69883947642649
$ws.Range("H6").Value = 13.16668431871898
$ws.Range("J6").Value = 9.909455443998008
$ws.Range("M6").Value = 17.28205949674647
$ws.Range("N6").Value = 17.66569508830805
$ws.Range("O6").Value = 19.69175130062039
$ws.Range("B7").Value = 7.661390209481102
$ws.Range("D7").Value = 10.29121105809111
$ws.Range("E7").Value = 14.11049442802857
$ws.Range("F7").Value = 27.07532645501286
$ws.Range("G7").Value = 25.80804434960562
$ws.Range("H7").Value = 13.16397935614277
$ws.Range("J7").Value = 9.896741289530288
$ws.Range("M7").Value = 17.43524503108592
$ws.Range("N7").Value = 17.69022513184995
$ws.Range("O7").Value = 19.70773468681212
$ws.Range("B8").Value = 7.753732210924177
$ws.Range("D8").Value = 10.24820356042513
$ws.Range("E8").Value = 14.01943506232824
$ws.Range("F8").Value = 27.19023543557205
$ws.Range("G8").Value = 26.30469564838683
$ws.Range("H8").Value = 13.1579367000847
$ws.Range("J8").Value = 9.843517668067813
$ws.Range("M8").Value = 18.09607392474
$ws.Range("N8").Value = 17.80337900698055
$ws.Range("O8").Value = 19.7896600131566
$ws.Range("B9").Value = 7.945188960833294
$ws.Range("D9").Value = 10.17632486226968
$ws.Range("E9").Value = 13.86105544047872
$ws.Range("F9").Value = 27.47477121148286
$ws.Range("G9").Value = 27.32248858833617
$ws.Range("H9").Value = 13.16723445767794
$ws.Range("J9").Value = 9.749452777294223
$ws.Range("M9").Value = 19.32682741395674
$ws.Range("N9").Value = 18.04170882069112
$ws.Range("O9").Value = 19.99044367601479
$ws.Range("B10").Value = 8.09089551023272
$ws.Range("D10").Value = 10.13111066183369
$ws.Range("E10").Value = 13.75697886212721
$ws.Range("F10").Value = 27.71863265417452
$ws.Range("G10").Value = 28.08865498578096
$ws.Range("H10").Value = 13.1869031909526
$ws.Range("J10").Value = 9.686579383322529
$ws.Range("M10").Value = 20.18295244814011
$ws.Range("N10").Value = 18.22570674493538
$ws.Range("O10").Value = 20.16154662464087
$ws.Range("B11").Value = 8.157998521071178
$ws.Range("D11").Value = 10.11219075718494
$ws.Range("E11").Value = 13.71229525649886
$ws.Range("F11").Value = 27.836809459689
$ws.Range("G11").Value = 28.43943799784849
$ws.Range("H11").Value = 13.19862025978119
$ws.Range("J11").Value = 9.659320830632215
$ws.Range("M11").Value = 20.56044965756783
$ws.Range("N11").Value = 18.31109589195117
$ws.Range("O11").Value = 20.24428903158983
$ws.Range("B12").Value = 8.183503321186334
$ws.Range("D12").Value = 10.10526333224627
$ws.Range("E12").Value = 13.69575728136191
$ws.Range("F12").Value = 27.88257104851826
$ws.Range("G12").Value = 28.5724504240593
$ws.Range("H12").Value = 13.20345332593907
$ws.Range("J12").Value = 9.649191112882011
$ws.Range("M12").Value = 20.70156685417991
$ws.Range("N12").Value = 18.34365175014733
$ws.Range("O12").Value = 20.27630642277817
$ws.Range("B13").Value = 8.178006584505543
$ws.Range("D13").Value = 10.10674472882818
$ws.Range("E13").Value = 13.69930200689793
$ws.Range("F13").Value = 27.87267103858672
$ws.Range("G13").Value = 28.54379814603918
$ws.Range("H13").Value = 13.20239486231673
$ws.Range("J13").Value = 9.651364176559001
$ws.Range("M13").Value = 20.67125793789399
$ws.Range("N13").Value = 18.33663079165235
$ws.Range("O13").Value = 20.26938080256099
$ws.Range("B14").Value = 8.160095060243547
$ws.Range("D14").Value = 10.11161608063097
$ws.Range("E14").Value = 13.71092699637013
$ws.Range("F14").Value = 27.84055422351205
$ws.Range("G14").Value = 28.45037824935382
$ws.Range("H14").Value = 13.19900995515151
$ws.Range("J14").Value = 9.658483597965532
$ws.Range("M14").Value = 20.57209679461275
$ws.Range("N14").Value = 18.31376998187628
$ws.Range("O14").Value = 20.24690952507408
$ws.Range("B15").Value = 8.14913532258004
$ws.Range("D15").Value = 10.11463080926502
$ws.Range("E15").Value = 13.71809748344538
$ws.Range("F15").Value = 27.82101240441542
$ws.Range("G15").Value = 28.3931748745878
$ws.Range("H15").Value = 13.19698810991731
$ws.Range("J15").Value = 9.66286950135642
$ws.Range("M15").Value = 20.51111587962234
$ws.Range("N15").Value = 18.29979518394373
$ws.Range("O15").Value = 20.23323373483489
$ws.Range("B16").Value = 8.086524656366096
$ws.Range("D16").Value = 10.13238029914323
$ws.Range("E16").Value = 13.75995259637301
$ws.Range("F16").Value = 27.71105267109342
$ws.Range("G16").Value = 28.06576328327429
$ws.Range("H16").Value = 13.18619300480938
$ws.Range("J16").Value = 9.688387764399113
$ws.Range("M16").Value = 20.15803154730908
$ws.Range("N16").Value = 18.22015846127148
$ws.Range("O16").Value = 20.15623618975233
$ws.Range("B17").Value = 8.048308313644778
$ws.Range("D17").Value = 10.14369127742575
$ws.Range("E17").Value = 13.78631100309932
$ws.Range("F17").Value = 27.64542967732035
$ws.Range("G17").Value = 27.86537881363133
$ws.Range("H17").Value = 13.1802784735342
$ws.Range("J17").Value = 9.704385903950623
$ws.Range("M17").Value = 19.93827873606756
$ws.Range("N17").Value = 18.17172068920671
$ws.Range("O17").Value = 20.11024295034062
$ws.Range("B18").Value = 8.026405910958072
$ws.Range("D18").Value = 10.15035218843506
$ws.Range("E18").Value = 13.80172215590989
$ws.Range("F18").Value = 27.60836879857887
$ws.Range("G18").Value = 27.75034322716106
$ws.Range("H18").Value = 13.17713743938237
$ws.Range("J18").Value = 9.713714042973903
$ws.Range("M18").Value = 19.81076317637119
$ws.Range("N18").Value = 18.14402073027613
$ws.Range("O18").Value = 20.08425204887743
$ws.Range("B19").Value = 8.019004359738583
$ws.Range("D19").Value = 10.15263409961214
$ws.Range("E19").Value = 13.80698312959773
$ws.Range("F19").Value = 27.5959389630104
$ws.Range("G19").Value = 27.7114364645148
$ws.Range("H19").Value = 13.17611880169214
$ws.Range("J19").Value = 9.716894123862676
$ws.Range("M19").Value = 19.76739995104356
$ws.Range("N19").Value = 18.134670169811
$ws.Range("O19").Value = 20.07553214726469
$ws.Range("B20").Value = 8.052368542529036
$ws.Range("D20").Value = 10.14247114771118
$ws.Range("E20").Value = 13.7834791803209
$ws.Range("F20").Value = 27.65234481109555
$ws.Range("G20").Value = 27.88668834912316
$ws.Range("H20").Value = 13.18088110359397
$ws.Range("J20").Value = 9.702669793783132
$ws.Range("M20").Value = 19.96178852379582
$ws.Range("N20").Value = 18.17686055344693
$ws.Range("O20").Value = 20.11509121807087
$ws.Range("B21").Value = 8.165353736062581
$ws.Range("D21").Value = 10.11017880991694
$ws.Range("E21").Value = 13.7075020673209
$ws.Range("F21").Value = 27.84996053967033
$ws.Range("G21").Value = 28.47781422917461
$ws.Range("H21").Value = 13.19999345460791
$ws.Range("J21").Value = 9.656387231515263
$ws.Range("M21").Value = 20.60127341398135
$ws.Range("N21").Value = 18.32047893741995
$ws.Range("O21").Value = 20.25349147648858
$ws.Range("B22").Value = 8.239735191649753
$ws.Range("D22").Value = 10.09045614101707
$ws.Range("E22").Value = 13.66007755800714
$ws.Range("F22").Value = 27.98498866765588
$ws.Range("G22").Value = 28.86512535649651
$ws.Range("H22").Value = 13.2147917926446
$ws.Range("J22").Value = 9.627260695920944
$ws.Range("M22").Value = 21.00848611508276
$ws.Range("N22").Value = 18.415618810908
$ws.Range("O22").Value = 20.34792403522127
$ws.Range("B23").Value = 8.199994839421903
$ws.Range("D23").Value = 10.10085598812266
$ws.Range("E23").Value = 13.68518475898333
$ws.Range("F23").Value = 27.91239479122993
$ws.Range("G23").Value = 28.6583676567582
$ws.Range("H23").Value = 13.20668334023168
$ws.Range("J23").Value = 9.642703626640415
$ws.Range("M23").Value = 20.79216511787402
$ws.Range("N23").Value = 18.36473129685766
$ws.Range("O23").Value = 20.29716668134994
$ws.Range("B24").Value = 8.050532697320891
$ws.Range("D24").Value = 10.14302227596898
$ws.Range("E24").Value = 13.78475864576084
$ws.Range("F24").Value = 27.64921640016553
$ws.Range("G24").Value = 27.87705377330274
$ws.Range("H24").Value = 13.18060784662069
$ws.Range("J24").Value = 9.703445240341079
$ws.Range("M24").Value = 19.95116340530512
$ws.Range("N24").Value = 18.17453635953123
$ws.Range("O24").Value = 20.11289790907812
$ws.Range("B25").Value = 7.89241521017636
$ws.Range("D25").Value = 10.194436166965
$ws.Range("E25").Value = 13.90174251277983
$ws.Range("F25").Value = 27.39158121593925
$ws.Range("G25").Value = 27.0433147530846
$ws.Range("H25").Value = 13.16245948089281
$ws.Range("J25").Value = 9.773801344838141
$ws.Range("M25").Value = 19.00176485642352
$ws.Range("N25").Value = 17.97558872400398
$ws.Range("O25").Value = 19.93191460752485
